$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the "From" hour value for the R30 rule (row 10, column C) from 18 to 1
$ws.Range("C10").Value = 1
